$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.434.57'
$ws.Range("E2").Value = '  +0.72%  '

$ws.Range("D3").Value = '1.838.61'
$ws.Range("E3").Value = '  -0.25%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9988'
$ws.Range("E4").Value = '  -1.16%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '243.32'
$ws.Range("E5").Value = '  -0.17%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.6276'
$ws.Range("E6").Value = '  +1.47%  '

$ws.Range("E7").Value = '  -1.04%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.07416'
$ws.Range("E8").Value = '  -0.36%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.2948'
$ws.Range("E9").Value = '  +0.02%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '23.43'
$ws.Range("E10").Value = '  +2.01%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07647'
$ws.Range("E11").Value = '  -0.91%  '

$ws.Range("D12").Value = '1.836.52'
$ws.Range("E12").Value = '  -0.29%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.014'
$ws.Range("E13").Value = '  +0.39%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.6759'
$ws.Range("E14").Value = '  +0.52%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '83.42'
$ws.Range("E15").Value = '  +0.64%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.000009352'
$ws.Range("E16").Value = '  +2.76%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '5.912'
$ws.Range("E17").Value = '  +0.58%  '

$ws.Range("D18").Value = '29.412.30'
$ws.Range("E18").Value = '  +0.66%  '

$ws.Range("D19").Value = '2.083.43'
$ws.Range("E19").Value = '  -0.28%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '237.56'
$ws.Range("E20").Value = '  +0.02%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '12.55'
$ws.Range("E21").Value = '  -0.49%  '

$ws.Range("E22").Value = '  -1.08%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.335'
$ws.Range("E23").Value = '  +2.30%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.0000'
$ws.Range("E24").Value = '  -1.44%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '159.00'
$ws.Range("E25").Value = '  -0.39%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.1414'
$ws.Range("E26").Value = '  -0.69%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.503'
$ws.Range("E27").Value = '  -0.10%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '17.76'
$ws.Range("E28").Value = '  -0.71%  '

$ws.Range("B29").Value = 'PancakeSwap'
$ws.Range("C29").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.496'
$ws.Range("E29").Value = '  -0.45%  '

$ws.Range("B30").Value = 'Hedera'
$ws.Range("C30").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.06009'
$ws.Range("E30").Value = '  +7.66%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.234'
$ws.Range("E31").Value = '  +1.17%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.094'
$ws.Range("E32").Value = '  -0.54%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.110'
$ws.Range("E33").Value = '  -0.81%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.867'
$ws.Range("E34").Value = '  +0.87%  '

$ws.Range("E35").Value = '  -0.01%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7258'
$ws.Range("E36").Value = '  -2.27%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.610'
$ws.Range("E37").Value = '  -1.73%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.882'
$ws.Range("E38").Value = '  +1.74%  '

$ws.Range("D39").Value = '1.217.59'
$ws.Range("E39").Value = '  +0.81%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.01758'
$ws.Range("E40").Value = '  -1.22%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '6.278'
$ws.Range("E41").Value = '  -2.56%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.9115'
$ws.Range("E42").Value = '  +0.03%  '

$ws.Range("E43").Value = '  -0.88%  '

$ws.Range("D44").Value = '1.997.20'
$ws.Range("E44").Value = '  +0.29%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '101.85'
$ws.Range("E45").Value = '  +0.26%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '65.62'
$ws.Range("E46").Value = '  +0.72%  '

$ws.Range("E47").Value = '  -1.50%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.5070'
$ws.Range("E48").Value = '  -1.32%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '9.220'
$ws.Range("E49").Value = '  +1.14%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.4061'
$ws.Range("E50").Value = '  +0.46%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.1140'
$ws.Range("E51").Value = '  +2.85%  '
